$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PUTWALL PICKING")
$ws.Range("A2:C8").ClearContents()
$ws.Cells.Item(2, 1).Value = 'ADOL798N.SEEMANNVAZQ'
$ws.Cells.Item(2, 2).Value = 24
$ws.Cells.Item(2, 3).Value = 11.52
$ws.Cells.Item(3, 1).Value = 'DIAN4065.ENTRIALGO'
$ws.Cells.Item(3, 2).Value = 37
$ws.Cells.Item(3, 3).Value = 17.76
$ws.Cells.Item(4, 1).Value = 'ZAKI0190.PHILLIPHORS'
$ws.Cells.Item(4, 2).Value = 2
$ws.Cells.Item(4, 3).Value = 0.96

$ws = $wb.Worksheets.Item("REGULAR PICK")
$ws.Range("A2:C2").ClearContents()
$ws.Cells.Item(2, 1).Value = 'ASHA1141.PAGE'
$ws.Cells.Item(2, 2).Value = 16
$ws.Cells.Item(2, 3).Value = 7.68
$ws.Cells.Item(3, 1).Value = 'DIAN4065.ENTRIALGO'
$ws.Cells.Item(3, 2).Value = 15
$ws.Cells.Item(3, 3).Value = 7.199999999999999
$ws.Cells.Item(4, 1).Value = 'KHINEHAYMAR.THAUNG'
$ws.Cells.Item(4, 2).Value = 1
$ws.Cells.Item(4, 3).Value = 0.48
$ws.Cells.Item(5, 1).Value = 'ZAKI0190.PHILLIPHORS'
$ws.Cells.Item(5, 2).Value = 1
$ws.Cells.Item(5, 3).Value = 0.48

$ws = $wb.Worksheets.Item("SINGLE PICK")
$ws.Range("A2:C10").ClearContents()
$ws.Cells.Item(2, 1).Value = 'BUDD0680.TENNAKOON'
$ws.Cells.Item(2, 2).Value = 5
$ws.Cells.Item(2, 3).Value = 2.4
$ws.Cells.Item(3, 1).Value = 'LOANA.MBONGO'
$ws.Cells.Item(3, 2).Value = 128
$ws.Cells.Item(3, 3).Value = 61.44
$ws.Cells.Item(4, 1).Value = 'MICA0432.RIZKALLAMAR'
$ws.Cells.Item(4, 2).Value = 5
$ws.Cells.Item(4, 3).Value = 2.4
$ws.Cells.Item(5, 1).Value = 'STAN9294.BAUER'
$ws.Cells.Item(5, 2).Value = 7
$ws.Cells.Item(5, 3).Value = 3.36
$ws.Cells.Item(6, 1).Value = 'WESL5337.CADETTE'
$ws.Cells.Item(6, 2).Value = 51
$ws.Cells.Item(6, 3).Value = 24.48
$ws.Cells.Item(7, 1).Value = 'WILDINE.JEUNE'
$ws.Cells.Item(7, 2).Value = 13
$ws.Cells.Item(7, 3).Value = 6.239999999999999

$ws = $wb.Worksheets.Item("REPLENISHMENT PICK")
$ws.Range("A2:C19").ClearContents()
$ws.Cells.Item(2, 1).Value = 'ADOL798N.SEEMANNVAZQ'
$ws.Cells.Item(2, 2).Value = 33
$ws.Cells.Item(2, 3).Value = 15.84
$ws.Cells.Item(3, 1).Value = 'ANJALI.BAKSHI'
$ws.Cells.Item(3, 2).Value = 3
$ws.Cells.Item(3, 3).Value = 1.44
$ws.Cells.Item(4, 1).Value = 'BOHD0676.KUSHLIAK'
$ws.Cells.Item(4, 2).Value = 6
$ws.Cells.Item(4, 3).Value = 2.88
$ws.Cells.Item(5, 1).Value = 'BUDD0680.TENNAKOON'
$ws.Cells.Item(5, 2).Value = 29
$ws.Cells.Item(5, 3).Value = 13.92
$ws.Cells.Item(6, 1).Value = 'DIAN4065.ENTRIALGO'
$ws.Cells.Item(6, 2).Value = 54
$ws.Cells.Item(6, 3).Value = 25.92
$ws.Cells.Item(7, 1).Value = 'IREN797N.CABRERA'
$ws.Cells.Item(7, 2).Value = 83
$ws.Cells.Item(7, 3).Value = 39.84
$ws.Cells.Item(8, 1).Value = 'JEEW9554.SITUMUDALIG'
$ws.Cells.Item(8, 2).Value = 48
$ws.Cells.Item(8, 3).Value = 23.04
$ws.Cells.Item(9, 1).Value = 'LOWRHY-OTIENO.JAOKO'
$ws.Cells.Item(9, 2).Value = 107
$ws.Cells.Item(9, 3).Value = 51.36
$ws.Cells.Item(10, 1).Value = 'MAKEDA.OLLIVIERRE'
$ws.Cells.Item(10, 2).Value = 96
$ws.Cells.Item(10, 3).Value = 46.08
$ws.Cells.Item(11, 1).Value = 'PATI2298.ATSIANGBE'
$ws.Cells.Item(11, 2).Value = 36
$ws.Cells.Item(11, 3).Value = 17.28
$ws.Cells.Item(12, 1).Value = 'PATR5027.AMEH'
$ws.Cells.Item(12, 2).Value = 43
$ws.Cells.Item(12, 3).Value = 20.64
$ws.Cells.Item(13, 1).Value = 'RAMI9087.SAIHI'
$ws.Cells.Item(13, 2).Value = 3
$ws.Cells.Item(13, 3).Value = 1.44
$ws.Cells.Item(14, 1).Value = 'RARG046N.YEBOAH'
$ws.Cells.Item(14, 2).Value = 73
$ws.Cells.Item(14, 3).Value = 35.04
$ws.Cells.Item(15, 1).Value = 'STAN9294.BAUER'
$ws.Cells.Item(15, 2).Value = 55
$ws.Cells.Item(15, 3).Value = 26.4
$ws.Cells.Item(16, 1).Value = 'THIE6554.DIALLO'
$ws.Cells.Item(16, 2).Value = 5
$ws.Cells.Item(16, 3).Value = 2.4
$ws.Cells.Item(17, 1).Value = 'WESL5337.CADETTE'
$ws.Cells.Item(17, 2).Value = 63
$ws.Cells.Item(17, 3).Value = 30.24
$ws.Cells.Item(18, 1).Value = 'WILDINE.JEUNE'
$ws.Cells.Item(18, 2).Value = 95
$ws.Cells.Item(18, 3).Value = 45.59999999999999
$ws.Cells.Item(19, 1).Value = 'ZAHIDGUL.MINHAS'
$ws.Cells.Item(19, 2).Value = 3
$ws.Cells.Item(19, 3).Value = 1.44

$ws = $wb.Worksheets.Item("QUICK MOVE")
$ws.Range("A2:C7").ClearContents()
$ws.Cells.Item(2, 1).Value = 'ADOL798N.SEEMANNVAZQ'
$ws.Cells.Item(2, 2).Value = 113
$ws.Cells.Item(2, 3).Value = 54.23999999999999
$ws.Cells.Item(3, 1).Value = 'DIAN4065.ENTRIALGO'
$ws.Cells.Item(3, 2).Value = 76
$ws.Cells.Item(3, 3).Value = 36.48
$ws.Cells.Item(4, 1).Value = 'ESSE0616.UDEH'
$ws.Cells.Item(4, 2).Value = 86
$ws.Cells.Item(4, 3).Value = 41.27999999999999
$ws.Cells.Item(5, 1).Value = 'MARI882N.ABDELKADER'
$ws.Cells.Item(5, 2).Value = 124
$ws.Cells.Item(5, 3).Value = 59.52
$ws.Cells.Item(6, 1).Value = 'RAMI9087.SAIHI'
$ws.Cells.Item(6, 2).Value = 8
$ws.Cells.Item(6, 3).Value = 3.84
$ws.Cells.Item(7, 1).Value = 'RARG046N.YEBOAH'
$ws.Cells.Item(7, 2).Value = 16
$ws.Cells.Item(7, 3).Value = 7.68
$ws.Cells.Item(8, 1).Value = 'STAN9294.BAUER'
$ws.Cells.Item(8, 2).Value = 25
$ws.Cells.Item(8, 3).Value = 12
$ws.Cells.Item(9, 1).Value = 'TANI2739.HOSSAINISLA'
$ws.Cells.Item(9, 2).Value = 20
$ws.Cells.Item(9, 3).Value = 9.6
$ws.Cells.Item(10, 1).Value = 'WESL5337.CADETTE'
$ws.Cells.Item(10, 2).Value = 63
$ws.Cells.Item(10, 3).Value = 30.24

$ws = $wb.Worksheets.Item("IDLE TIME")
$ws.Range("A2:B22").ClearContents()
$ws.Cells.Item(2, 1).Value = 'ADOL798N.SEEMANNVAZQ'
$ws.Cells.Item(2, 2).Value = 52
$ws.Cells.Item(3, 1).Value = 'ANJALI.BAKSHI'
$ws.Cells.Item(3, 2).Value = 124
$ws.Cells.Item(4, 1).Value = 'ASHA1141.PAGE'
$ws.Cells.Item(4, 2).Value = 25
$ws.Cells.Item(5, 1).Value = 'BOHD0676.KUSHLIAK'
$ws.Cells.Item(5, 2).Value = 119
$ws.Cells.Item(6, 1).Value = 'BUDD0680.TENNAKOON'
$ws.Cells.Item(6, 2).Value = 92
$ws.Cells.Item(7, 1).Value = 'DIAN4065.ENTRIALGO'
$ws.Cells.Item(7, 2).Value = 39
$ws.Cells.Item(8, 1).Value = 'ESSE0616.UDEH'
$ws.Cells.Item(8, 2).Value = 97
$ws.Cells.Item(9, 1).Value = 'IREN797N.CABRERA'
$ws.Cells.Item(9, 2).Value = 28
$ws.Cells.Item(10, 1).Value = 'JEEW9554.SITUMUDALIG'
$ws.Cells.Item(10, 2).Value = 22
$ws.Cells.Item(11, 1).Value = 'KHINEHAYMAR.THAUNG'
$ws.Cells.Item(11, 2).Value = 88
$ws.Cells.Item(12, 1).Value = 'LOANA.MBONGO'
$ws.Cells.Item(12, 2).Value = 33
$ws.Cells.Item(13, 1).Value = 'LOWRHY-OTIENO.JAOKO'
$ws.Cells.Item(13, 2).Value = 21
$ws.Cells.Item(14, 1).Value = 'MAKEDA.OLLIVIERRE'
$ws.Cells.Item(14, 2).Value = 37
$ws.Cells.Item(15, 1).Value = 'MARI882N.ABDELKADER'
$ws.Cells.Item(15, 2).Value = 85
$ws.Cells.Item(16, 1).Value = 'MICA0432.RIZKALLAMAR'
$ws.Cells.Item(16, 2).Value = 124
$ws.Cells.Item(17, 1).Value = 'PATI2298.ATSIANGBE'
$ws.Cells.Item(17, 2).Value = 74
$ws.Cells.Item(18, 1).Value = 'PATR5027.AMEH'
$ws.Cells.Item(18, 2).Value = 49
$ws.Cells.Item(19, 1).Value = 'RAMI9087.SAIHI'
$ws.Cells.Item(19, 2).Value = 117
$ws.Cells.Item(20, 1).Value = 'RARG046N.YEBOAH'
$ws.Cells.Item(20, 2).Value = 41
$ws.Cells.Item(21, 1).Value = 'STAN9294.BAUER'
$ws.Cells.Item(21, 2).Value = 21
$ws.Cells.Item(22, 1).Value = 'TANI2739.HOSSAINISLA'
$ws.Cells.Item(22, 2).Value = 116
$ws.Cells.Item(23, 1).Value = 'THIE6554.DIALLO'
$ws.Cells.Item(23, 2).Value = 120
$ws.Cells.Item(24, 1).Value = 'WESL5337.CADETTE'
$ws.Cells.Item(24, 2).Value = 28
$ws.Cells.Item(25, 1).Value = 'WILDINE.JEUNE'
$ws.Cells.Item(25, 2).Value = 36
$ws.Cells.Item(26, 1).Value = 'ZAHIDGUL.MINHAS'
$ws.Cells.Item(26, 2).Value = 121
$ws.Cells.Item(27, 1).Value = 'ZAKI0190.PHILLIPHORS'
$ws.Cells.Item(27, 2).Value = 112

$ws = $wb.Worksheets.Item("Total Units picked by hour")
$ws.Cells.Item(2, 1).Value = 20
$ws.Cells.Item(2, 2).Value = -16
$ws.Cells.Item(2, 3).Value = -40
$ws.Cells.Item(2, 4).Value = -396
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(3, 1).Value = 21
$ws.Cells.Item(3, 2).Value = -1
$ws.Cells.Item(3, 3).Value = -137
$ws.Cells.Item(3, 4).Value = -399
$ws.Cells.Item(3, 5).Value = -63
$ws.Cells.Item(4, 1).Value = 22
$ws.Cells.Item(4, 2).Value = -16
$ws.Cells.Item(4, 3).Value = -32
$ws.Cells.Item(4, 4).Value = -40
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(5, 1).Value = 23
$ws.Cells.Item(5, 2).Value = 0
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(6, 1).Value = 'Total'
$ws.Cells.Item(6, 2).Value = -33
$ws.Cells.Item(6, 3).Value = -209
$ws.Cells.Item(6, 4).Value = -835
$ws.Cells.Item(6, 5).Value = -63
